# Hide filter in locations lists:
# Reorder/rename the header row of the terminal import sample sheet so
# that Terminal ID / Merchant ID come first, and drop the "...Id" style
# labels in favor of friendlier, UI-facing names.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row values (A1:J1)
$ws.Range("A1").Value = "Terminal ID"
$ws.Range("B1").Value = "Merchant ID"
$ws.Range("C1").Value = "ErrandChannel"
$ws.Range("D1").Value = "POSType"
$ws.Range("E1").Value = "PhoneNumber"
$ws.Range("F1").Value = "Latitude"
$ws.Range("G1").Value = "Longitude"
$ws.Range("H1").Value = "Zone"
$ws.Range("I1").Value = "Address"
$ws.Range("J1").Value = "Landmark"

# Resize the columns to (re-)fit their new header text.
$ws.Columns.Item(1).ColumnWidth = 9.5
$ws.Columns.Item(2).ColumnWidth = 10.333333333333334
$ws.Columns.Item(3).ColumnWidth = 13.833333333333334
$ws.Columns.Item(4).ColumnWidth = 9.333333333333334
$ws.Columns.Item(6).ColumnWidth = 6.333333333333333
$ws.Columns.Item(7).ColumnWidth = 8.0
$ws.Columns.Item(8).ColumnWidth = 5.333333333333333
$ws.Columns.Item(9).ColumnWidth = 6.833333333333333
$ws.Columns.Item(10).ColumnWidth = 8.166666666666666

# Move the active selection back to the top of the sheet.
$ws.Range("D1").Select()
